# Add a new data row (row 3) to the worksheet, mirroring the structure of
# the existing row 2, and extend the sheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 112468257
$ws.Range("B3").Value = 98267
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 935
$ws.Range("F3").Value = "Vildris"
$ws.Range("G3").Value = "Leersia oryzoides"
$ws.Range("H3").Value = "(L.) Sw."

# "Antal" (I3) must stay textual ("100"), not become a number. Force text
# format before writing, then drop the extra formatting so the cell keeps
# the workbook's default style (matching the other plain-text cells).
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "100"
$ws.Range("I3").ClearFormats()

$ws.Range("J3").Value = "stjälkar/strån/skott"
$ws.Range("K3").Value = "fullt utvecklade blad"

# L3 / N3 (Kön / Metod) are present but blank in the source data, exactly
# like L2 / N2 on the row above — copy those empty cells down so row 3 gets
# an actual (empty) cell rather than no cell at all.
$ws.Range("L2").Copy($ws.Range("L3"))
$ws.Range("N2").Copy($ws.Range("N3"))

$ws.Range("P3").Value = "Hjörnereds sjösystem, vid nordvästligaste vägbanken, Hl"
$ws.Range("Q3").Value = 388943
$ws.Range("R3").Value = 6264177
$ws.Range("S3").Value = 5
$ws.Range("T3").Value = "Halland"
$ws.Range("U3").Value = "Laholm"
$ws.Range("V3").Value = "Halland"
$ws.Range("W3").Value = "Ysby"
$ws.Range("X3").Value = "N-Lah-1477"

# Startdatum / Slutdatum (Y3 / AA3) are textual dates ("2023-09-27"), not
# Excel date serials — same text-forcing trick as I3.
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-09-27"
$ws.Range("Y3").ClearFormats()

$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-09-27"
$ws.Range("AA3").ClearFormats()

$ws.Range("AD3").Value = $False
$ws.Range("AE3").Value = $False

# AF3 (Bestämningsmetod) blank like AF2.
$ws.Range("AF2").Copy($ws.Range("AF3"))

$ws.Range("AG3").Value = $False

# AT3 (Bestämningsår) blank like AT2.
$ws.Range("AT2").Copy($ws.Range("AT3"))

$ws.Range("AW3").Value = "Lars-Erik Magnusson"
$ws.Range("AX3").Value = "Lars-Erik Magnusson"
$ws.Range("AY3").Value = "Floraväkteri Sverige"
